$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> [new C value, new E value]
$updates = @{
    3   = @(249336, 1036489581)
    53  = @(141688, 590078732)
    91  = @(151249, 483346675)
    92  = @(409341, 1597795335)
    93  = @(209683, 1310346798)
    94  = @(94246, 919357089)
    95  = @(50813, 935119531)
    96  = @(17331, 797952512)
    104 = @(135345, 273027754)
    116 = @(4567, 20680072)
    118 = @(981, 11896140)
    145 = @(11834, 182736103)
    163 = @(70985, 131769402)
    173 = @(96866, 327978644)
    174 = @(226115, 900825976)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

$wb.Save()
